$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 552
$ws.Range("F4").Value = 6020
$ws.Range("F5").Value = 76
$ws.Range("F6").Value = 64
$ws.Range("F12").Value = 694
$ws.Range("F13").Value = 1621
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 1660
$ws.Range("F16").Value = 574
$ws.Range("F17").Value = 201
$ws.Range("F18").Value = 657
$ws.Range("F19").Value = 4695
$ws.Range("F20").Value = 115
$ws.Range("F21").Value = 52
$ws.Range("F23").Value = 3364
$ws.Range("F24").Value = 827
$ws.Range("F25").Value = 26
$ws.Range("F27").Value = 19
$ws.Range("F28").Value = 2357
$ws.Range("F30").Value = 344
$ws.Range("F31").Value = 5
$ws.Range("F36").Value = 27
$ws.Range("F39").Value = 1287
$ws.Range("F40").Value = 1273

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 16
$ws.Range("F15").Value = 70
$ws.Range("F22").Value = 243

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 767
$ws.Range("F4").Value = 214

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 552
$ws.Range("F6").Value = 767
$ws.Range("F7").Value = 214
$ws.Range("F8").Value = 6020
$ws.Range("F19").Value = 16
$ws.Range("F22").Value = 1621
$ws.Range("F24").Value = 1660
$ws.Range("F25").Value = 574
$ws.Range("F26").Value = 201
$ws.Range("F27").Value = 657
$ws.Range("F28").Value = 4695
$ws.Range("F29").Value = 52
$ws.Range("F30").Value = 679
$ws.Range("F31").Value = 3364
$ws.Range("F32").Value = 827
$ws.Range("F35").Value = 19
$ws.Range("F36").Value = 2357
$ws.Range("F38").Value = 344
$ws.Range("F42").Value = 243
$ws.Range("F45").Value = 27
